$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores dot-grouped numbers as plain text.
# Values that look like a plain decimal (single dot, e.g. "581.84") would
# otherwise be auto-coerced to a floating point number by COM, so for those
# we briefly force Text format, assign, then restore the default style so the
# cell formatting stays exactly as it was (no "s" attribute).
$ws.Range("D2").Value = "63.514.48"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").Value = "3.475.67"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("D7").Value = "3.476.46"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.125"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.405"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.26%  "
$ws.Range("D13").Value = "4.074.32"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.06%  "
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "3.470.66"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "63.479.18"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.565"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "3.622.49"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("E28").Value = "  -3.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "169.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").Value = "3.513.35"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0765"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.800"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "2.618.95"
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.76%  "
